# Fix: Update dashboard and KPIs
# Refresh the computed mtm/faceValue/profitLoss (and a few driver volume)
# figures on Sheet1 for the m0/m1/m2/m3/a0 maturities (rows 2-17) to match
# the latest upstream net_products recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 4335656.957893601
$ws.Range("N2").Value = 4867736.9578936
$ws.Range("Q2").Value = -8342153.656287315
$ws.Range("E3").Value = 1.915432796301079
$ws.Range("G3").Value = -2.134567203698921
$ws.Range("H3").Value = -9.5
$ws.Range("L3").Value = -9341285.450362897
$ws.Range("M3").Value = -168778.0958781447
$ws.Range("N3").Value = 387622.3041218542
$ws.Range("Q3").Value = -9510063.546241041
$ws.Range("E4").Value = -9.232429167111111
$ws.Range("G4").Value = -10.14242916711111
$ws.Range("H4").Value = -15.5
$ws.Range("L4").Value = -7353661.083225883
$ws.Range("M4").Value = -2619207.28809984
$ws.Range("N4").Value = -2073972.88809984
$ws.Range("Q4").Value = -9972868.371325726
$ws.Range("M5").Value = -5738073.145677392
$ws.Range("N5").Value = -5171372.422797393
$ws.Range("Q5").Value = -9615406.939540802
$ws.Range("L6").Value = 1318935.414203487
$ws.Range("M6").Value = -5635800
$ws.Range("N6").Value = -3989700
$ws.Range("Q6").Value = -4316864.585796515
$ws.Range("L7").Value = 1447933.303151536
$ws.Range("M7").Value = -5347036.8
$ws.Range("N7").Value = -3904320
$ws.Range("Q7").Value = -3899103.496848463
$ws.Range("L8").Value = 1248627.414203487
$ws.Range("M8").Value = -5565492
$ws.Range("N8").Value = -3892980
$ws.Range("Q8").Value = -4316864.585796515
$ws.Range("L9").Value = 951669.1105195042
$ws.Range("M9").Value = -5129280
$ws.Range("N9").Value = -3493800
$ws.Range("Q9").Value = -4177610.889480497
$ws.Range("L10").Value = 997229.8142034872
$ws.Range("M10").Value = -5314094.4
$ws.Range("N10").Value = -3688380
$ws.Range("Q10").Value = -4316864.585796515
$ws.Range("L11").Value = 962973.1105195042
$ws.Range("M11").Value = -5140584
$ws.Range("N11").Value = -3600000
$ws.Range("Q11").Value = -4177610.889480497
$ws.Range("L12").Value = 1075126.614203486
$ws.Range("M12").Value = -5391991.2
$ws.Range("N12").Value = -3827880
$ws.Range("Q12").Value = -4316864.585796515
$ws.Range("L13").Value = 1300335.414203487
$ws.Range("M13").Value = -5617200
$ws.Range("N13").Value = -4093860
$ws.Range("Q13").Value = -4316864.585796515
$ws.Range("L14").Value = 1457829.110519504
$ws.Range("M14").Value = -5635440
$ws.Range("N14").Value = -4163400
$ws.Range("Q14").Value = -4177610.889480497
$ws.Range("L15").Value = 1516021.014203486
$ws.Range("M15").Value = -5832885.6
$ws.Range("N15").Value = -4278000
$ws.Range("Q15").Value = -4316864.585796515
$ws.Range("L16").Value = 1256013.110519504
$ws.Range("M16").Value = -5433624
$ws.Range("N16").Value = -3866400
$ws.Range("Q16").Value = -4177610.889480497
$ws.Range("L17").Value = 1358813.814203487
$ws.Range("M17").Value = -5675678.4
$ws.Range("N17").Value = -4064100
$ws.Range("Q17").Value = -4316864.585796515
